$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26+ down by one
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with its data
$ws.Range("A26").Value = 11
$ws.Range("B26").Value = "Vega Monumental Concepción"
$ws.Range("C26").Value = "Bíobío"
$ws.Range("D26").Value = 44811
$ws.Range("E26").Value = 8
$ws.Range("F26").Value = 100112001
$ws.Range("G26").Value = "Berenjena"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 12000
$ws.Range("L26").Value = 13000
$ws.Range("M26").Value = 12500
$ws.Range("N26").Value = "$/caja 60 unidades"
$ws.Range("O26").Value = "Región de Arica y Parinacota"
$ws.Range("P26").Value = 208
$ws.Range("Q26").Value = 60
$ws.Range("R26").Value = "Hortaliza"
